$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove duplicate id_heating_technology rows (second occurrences of 29, 210, 211, 33)
# Delete from bottom to top so row numbers of rows above stay stable.
$ws.Rows("22:22").Delete()
$ws.Rows("16:16").Delete()
$ws.Rows("14:14").Delete()
$ws.Rows("12:12").Delete()

$ws.Range("B9").Select()
